# Weekly update: a new data point is inserted at row 387 ("Fruta, Macroferia
# Regional de Talca - Piña"), pushing every following row's data down by one
# row. The last existing row (438) ends up duplicated into the newly
# appended row 439. Columns A,B,C,E,F,G,H,I,J,K,R are constant across this
# block of rows, so only D,L,M,N,O,P,Q,S,T need to move.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 387
$lastRow  = 438
$newLastRow = $lastRow + 1

$cols = @("D","L","M","N","O","P","Q","S","T")

# Shift rows down: process from bottom to top so we don't clobber data
# before it has been read.
for ($r = $lastRow; $r -ge $firstRow; $r--) {
    foreach ($col in $cols) {
        $src = $ws.Range("$col$r")
        $dst = $ws.Range("$col$($r + 1)")
        $dst.Value = $src.Value2
    }
}

# Row 387 now gets the brand-new data point for the latest week.
$ws.Range("D387").Value = 45154
$ws.Range("L387").Value = "Segunda"
$ws.Range("M387").Value = 250
$ws.Range("N387").Value = 21000
$ws.Range("O387").Value = 21000
$ws.Range("P387").Value = 21000
$ws.Range("Q387").Value = "$/caja 14 unidades"
$ws.Range("S387").Value = 1500
$ws.Range("T387").Value = 14

# The constant columns for the new row (439) must be filled in too, since
# that row didn't exist before.
$ws.Range("A439").Value = $ws.Range("A438").Value2
$ws.Range("B439").Value = $ws.Range("B438").Value2
$ws.Range("C439").Value = $ws.Range("C438").Value2
$ws.Range("E439").Value = $ws.Range("E438").Value2
$ws.Range("F439").Value = $ws.Range("F438").Value2
$ws.Range("G439").Value = $ws.Range("G438").Value2
$ws.Range("H439").Value = $ws.Range("H438").Value2
$ws.Range("I439").Value = $ws.Range("I438").Value2
$ws.Range("J439").Value = $ws.Range("J438").Value2
$ws.Range("K439").Value = $ws.Range("K438").Value2
$ws.Range("R439").Value = $ws.Range("R438").Value2

# Column D carries a date number format (style index 2 in the original
# file); new cells default to General, so copy the format explicitly.
$ws.Range("D439").NumberFormat = $ws.Range("D438").NumberFormat
